$d = $word.ActiveDocument

# 1. Fill in the previously-empty bullet under "What are some limitations of this dataset?"
#    with the new observation about the 33%/50% success rate discrepancy.
$p6 = $d.Paragraphs(6)
$p6.Range.InsertBefore("The background noted that 33% of projects succeeded, however, the sample presented has over 50% successful projects. Thus, the data is not representative of the actual success and failure rate which can skew the analysis on the trend for success. ")

# 2. Expand "the days project was live" -> "the account of days the project was live"
$d.Content.Find.Execute("We could create a table that filters by the days project was live.", $true, $false, $false, $false, $false, $true, 1, $false, "We could create a table that filters by the account of days the project was live.", 2)

# 3. Expand "the significance of goal amount" -> "the significance of the goal amount"
$d.Content.Find.Execute("We could create a table that observes the significance of goal amount.", $true, $false, $false, $false, $false, $true, 1, $false, "We could create a table that observes the significance of the goal amount.", 2)
